# Populate the data grid using the modern Range-based cell access
# (replaces the old/deprecated Cells.Item(...) access pattern).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: column headers
$ws.Range('A1').Value = 'NAME'
$ws.Range('B1').Value = 'DESCRIPTION'
$ws.Range('C1').Value = 'NAME'
$ws.Range('D1').Value = 'PROJECT_ID'
$ws.Range('E1').Value = 'SHEET'
$ws.Range('F1').Value = 'CONTENT_TYPE_ID'
$ws.Range('G1').Value = 'OBJECT_ID'
$ws.Range('H1').Value = 'TYPE'
$ws.Range('I1').Value = 'TEMPLATE_ID'
$ws.Range('J1').Value = 'DB_FIELD'
$ws.Range('K1').Value = 'DATA_TYPE'
$ws.Range('L1').Value = 'CONFIGURABLE_SETTING'
$ws.Range('M1').Value = 'LABEL'
$ws.Range('N1').Value = 'EXCEPTIONS'
$ws.Range('O1').Value = 'REQUIRED'
$ws.Range('P1').Value = 'ACTIVE'
$ws.Range('Q1').Value = 'DISABLED'

# Row 2: data record
$ws.Range('A2').Value = 'org3'
$ws.Range('B2').Value = 'desc3'
$ws.Range('C2').Value = 'temp2'
$ws.Range('D2').Value = 1
$ws.Range('E2').Value = '[{"_id":"5b20b297fc28b53b80bb588b","index":0,"guid":"b9cc2193-c475-4f86-9c5d-ba1c1eb17cfa","isActive":false,"balance":"$2,544.76","picture":"http://placehold.it/32x32","age":35,"eyeColor":"blue","name":"Deana Cohen"},{"_id":"5b20b2973f8a4a90ad8c35fc","index":1,"guid":"b2b0f4d5-39e1-4321-b300-3c08c0eacf4b","isActive":false,"balance":"$2,252.69","picture":"http://placehold.it/32x32","age":25,"eyeColor":"brown","name":"Chen Hooper"},{"_id":"5b20b2977df4355aa7b993af","index":2,"guid":"1f0c0829-f9ab-4eb6-b106-377e41086ddf","isActive":false,"balance":"$1,116.92","picture":"http://placehold.it/32x32","age":37,"eyeColor":"brown","name":"Arnold Wilson"},{"_id":"5b20b2973400739f91af8f3f","index":3,"guid":"14662d70-6e10-4d6b-bba8-24322fe79ad3","isActive":false,"balance":"$1,231.11","picture":"http://placehold.it/32x32","age":23,"eyeColor":"green","name":"Berta Clark"},{"_id":"5b20b297b9504bbe4b26303f","index":4,"guid":"966ec57f-f3ce-4eab-9d2f-636f3f1f89f1","isActive":true,"balance":"$1,310.59","picture":"http://placehold.it/32x32","age":35,"eyeColor":"brown","name":"Vaughn Ellison"},{"_id":"5b20b297f043757eca7cdedc","index":5,"guid":"22469496-b935-448a-b206-3fa34aa4341c","isActive":false,"balance":"$3,774.91","picture":"http://placehold.it/32x32","age":28,"eyeColor":"green","name":"Jacqueline Nixon"}]'
$ws.Range('G2').Value = 50
$ws.Range('H2').Value = 'jojojojo'
$ws.Range('I2').Value = 1
$ws.Range('J2').Value = 'conftemp2'
$ws.Range('K2').Value = 'int'
$ws.Range('L2').Value = '[{"_id":"5b20b333866cdbec105144f0","index":0,"guid":"fe4c40ba-1150-4ec7-be2b-68364b0ce1c3","isActive":true,"balance":"$1,557.90","picture":"http://placehold.it/32x32","age":29,"eyeColor":"brown","name":"Priscilla Townsend"},{"_id":"5b20b333f2b68773771fe88d","index":1,"guid":"ba7f3edf-f3db-46cf-9b68-cad0f26a159b","isActive":true,"balance":"$1,287.93","picture":"http://placehold.it/32x32","age":29,"eyeColor":"green","name":"Peggy Cherry"},{"_id":"5b20b333545e5f1a321976f6","index":2,"guid":"6f58410c-adb4-4ecd-87b3-59fd1dcbe936","isActive":false,"balance":"$2,818.10","picture":"http://placehold.it/32x32","age":38,"eyeColor":"brown","name":"Shepherd Barry"},{"_id":"5b20b333a66072a9828aba04","index":3,"guid":"774ddfc3-f81c-4e8d-a1bd-9fc3a1dbf6f9","isActive":true,"balance":"$3,666.68","picture":"http://placehold.it/32x32","age":24,"eyeColor":"blue","name":"Patton Hansen"},{"_id":"5b20b333a8ac10fa3032837c","index":4,"guid":"1f36bd52-18d8-46a7-a222-579ee34f0047","isActive":false,"balance":"$3,146.98","picture":"http://placehold.it/32x32","age":39,"eyeColor":"brown","name":"House Ross"}]'
$ws.Range('M2').Value = 'asdsoimas'
$ws.Range('N2').Value = '[{"_id":"5b20b33d44b7ab485e26c63c","index":0,"guid":"70cd3fa3-715c-4b9c-adf0-055ac8224fb6","isActive":false,"balance":"$3,620.65","picture":"http://placehold.it/32x32","age":35,"eyeColor":"blue","name":"Vance Daniel"},{"_id":"5b20b33d5b49487ec4c8b9a8","index":1,"guid":"44464f93-900c-4baa-9be9-65bf83b34ff1","isActive":true,"balance":"$3,784.97","picture":"http://placehold.it/32x32","age":31,"eyeColor":"green","name":"Tania Preston"},{"_id":"5b20b33d2badf4414a697d92","index":2,"guid":"b4273a3e-0ed4-4522-a181-499307ce3169","isActive":true,"balance":"$1,365.76","picture":"http://placehold.it/32x32","age":24,"eyeColor":"brown","name":"Lea Whitehead"},{"_id":"5b20b33d8daddacd45035873","index":3,"guid":"422625a5-6f43-4701-ad85-308e315b1e2f","isActive":false,"balance":"$2,686.75","picture":"http://placehold.it/32x32","age":27,"eyeColor":"blue","name":"Rosemarie Duffy"},{"_id":"5b20b33d2ae9597c7b72050b","index":4,"guid":"9879a309-cae6-4df6-b174-a091bbcfe1f4","isActive":true,"balance":"$2,054.49","picture":"http://placehold.it/32x32","age":29,"eyeColor":"brown","name":"Richmond Kelly"},{"_id":"5b20b33df2a180bf70f6c60b","index":5,"guid":"15bda1b1-3664-4e56-8d5a-7e57fd09b28c","isActive":false,"balance":"$3,488.15","picture":"http://placehold.it/32x32","age":22,"eyeColor":"green","name":"Ware Farrell"}]'
$ws.Range('O2').Value = $false
$ws.Range('P2').Value = $true
$ws.Range('Q2').Value = $true

# Row 3: data record
$ws.Range('A3').Value = 'org3'
$ws.Range('B3').Value = 'desc3'
$ws.Range('C3').Value = 'temp2'
$ws.Range('D3').Value = 1
$ws.Range('E3').Value = '[{"_id":"5b20b297fc28b53b80bb588b","index":0,"guid":"b9cc2193-c475-4f86-9c5d-ba1c1eb17cfa","isActive":false,"balance":"$2,544.76","picture":"http://placehold.it/32x32","age":35,"eyeColor":"blue","name":"Deana Cohen"},{"_id":"5b20b2973f8a4a90ad8c35fc","index":1,"guid":"b2b0f4d5-39e1-4321-b300-3c08c0eacf4b","isActive":false,"balance":"$2,252.69","picture":"http://placehold.it/32x32","age":25,"eyeColor":"brown","name":"Chen Hooper"},{"_id":"5b20b2977df4355aa7b993af","index":2,"guid":"1f0c0829-f9ab-4eb6-b106-377e41086ddf","isActive":false,"balance":"$1,116.92","picture":"http://placehold.it/32x32","age":37,"eyeColor":"brown","name":"Arnold Wilson"},{"_id":"5b20b2973400739f91af8f3f","index":3,"guid":"14662d70-6e10-4d6b-bba8-24322fe79ad3","isActive":false,"balance":"$1,231.11","picture":"http://placehold.it/32x32","age":23,"eyeColor":"green","name":"Berta Clark"},{"_id":"5b20b297b9504bbe4b26303f","index":4,"guid":"966ec57f-f3ce-4eab-9d2f-636f3f1f89f1","isActive":true,"balance":"$1,310.59","picture":"http://placehold.it/32x32","age":35,"eyeColor":"brown","name":"Vaughn Ellison"},{"_id":"5b20b297f043757eca7cdedc","index":5,"guid":"22469496-b935-448a-b206-3fa34aa4341c","isActive":false,"balance":"$3,774.91","picture":"http://placehold.it/32x32","age":28,"eyeColor":"green","name":"Jacqueline Nixon"}]'
$ws.Range('G3').Value = 50
$ws.Range('H3').Value = 'jojojojo'
$ws.Range('I3').Value = 1
$ws.Range('J3').Value = 'conftemp2'
$ws.Range('K3').Value = 'table'
$ws.Range('L3').Value = '[{"name":"Rodgers Delacruz","gender":"male","greeting":"Hello, Rodgers Delacruz! You have 10 unread messages.","favoriteFruit":"banana"},{"name":"Garcia Moon","gender":"male","greeting":"Hello, Garcia Moon! You have 5 unread messages.","favoriteFruit":"strawberry"},{"name":"Lane Patrick","gender":"male","greeting":"Hello, Lane Patrick! You have 3 unread messages.","favoriteFruit":"strawberry"},{"name":"Moody Barrera","gender":"male","greeting":"Hello, Moody Barrera! You have 7 unread messages.","favoriteFruit":"strawberry"},{"name":"Meagan Hester","gender":"female","greeting":"Hello, Meagan Hester! You have 3 unread messages.","favoriteFruit":"strawberry"},{"name":"Lewis Morse","gender":"male","greeting":"Hello, Lewis Morse! You have 3 unread messages.","favoriteFruit":"banana"},{"name":"Bonnie Cross","gender":"female","greeting":"Hello, Bonnie Cross! You have 1 unread messages.","favoriteFruit":"banana"}]'
$ws.Range('M3').Value = 'asdkajsdasn'
$ws.Range('N3').Value = '[{"name":"Rodgers Delacruz","gender":"male","greeting":"Hello, Rodgers Delacruz! You have 10 unread messages.","favoriteFruit":"banana"},{"name":"Garcia Moon","gender":"male","greeting":"Hello, Garcia Moon! You have 5 unread messages.","favoriteFruit":"strawberry"},{"name":"Lane Patrick","gender":"male","greeting":"Hello, Lane Patrick! You have 3 unread messages.","favoriteFruit":"strawberry"},{"name":"Moody Barrera","gender":"male","greeting":"Hello, Moody Barrera! You have 7 unread messages.","favoriteFruit":"strawberry"},{"name":"Meagan Hester","gender":"female","greeting":"Hello, Meagan Hester! You have 3 unread messages.","favoriteFruit":"strawberry"},{"name":"Lewis Morse","gender":"male","greeting":"Hello, Lewis Morse! You have 3 unread messages.","favoriteFruit":"banana"},{"name":"Bonnie Cross","gender":"female","greeting":"Hello, Bonnie Cross! You have 1 unread messages.","favoriteFruit":"banana"}]'
$ws.Range('O3').Value = $false
$ws.Range('P3').Value = $true
$ws.Range('Q3').Value = $true
